$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Incorporate updated data from upstream processes through 2024:
# the "Open year" = 2024 row (row 26) gets revised facility counts for
# Energy Storage (column C) and Solar (column E).
$ws.Range("C26").Value = 1
$ws.Range("E26").Value = 8
